$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.059.25"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.651.00"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.97"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5301"
$ws.Range("E6").Value = "  +1.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("E8").Value = "  -2.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06300"
$ws.Range("E9").Value = "  -0.44%  "
$ws.Range("E10").Value = "  -3.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07736"
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.476"
$ws.Range("E12").Value = "  +0.83%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.600.24"
$ws.Range("E13").Value = "  -3.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5446"
$ws.Range("E14").Value = "  -0.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0₅8110"
$ws.Range("E15").Value = "  -1.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.13"
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.084.28"
$ws.Range("E17").Value = "  -0.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.003"
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.555"
$ws.Range("E19").Value = "  -2.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.56"
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("E21").Value = "  -1.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.989"
$ws.Range("E22").Value = "  -1.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.004"
$ws.Range("E23").Value = "  -0.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "139.33"
$ws.Range("E24").Value = "  +0.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1244"
$ws.Range("E25").Value = "  -0.27%  "
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.20"
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.438"
$ws.Range("E28").Value = "  +0.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05919"
$ws.Range("E29").Value = "  -1.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.280"
$ws.Range("E30").Value = "  -0.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.504"
$ws.Range("E31").Value = "  -1.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.241"
$ws.Range("E32").Value = "  -3.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.549"
$ws.Range("E33").Value = "  -6.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.415"
$ws.Range("E34").Value = "  +0.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9414"
$ws.Range("E35").Value = "  -4.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.750"
$ws.Range("E36").Value = "  -0.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5666"
$ws.Range("E37").Value = "  -4.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01606"
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.833"
$ws.Range("E39").Value = "  -2.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8447"
$ws.Range("E40").Value = "  -2.26%  "
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.008.76"
$ws.Range("E42").Value = "  -3.08%  "
$ws.Range("E43").Value = "  +0.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.800.58"
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "56.90"
$ws.Range("E45").Value = "  -0.60%  "
$ws.Range("E46").Value = "  -2.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.008"
$ws.Range("E47").Value = "  +0.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4296"
$ws.Range("E48").Value = "  +1.58%  "
$ws.Range("E49").Value = "  +0.63%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.874"
$ws.Range("E50").Value = "  -2.95%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05152"
$ws.Range("E51").Value = "  -0.53%  "
